# Remove rpart and add 50percent cutoff based on JKL feedback.
#
# The "MethylScoreAML Categorical" (column C) / "MethylScoreAML_cat_bin"
# (column D) labels were previously derived from an rpart split. They are
# replaced here with a simple 50% cutoff on the numeric MethylScoreAML
# score in column B: the top half of samples (by score) are labelled
# "High" (bin = 1) and the bottom half are labelled "Low" (bin = 0).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data starts on row 2 (row 1 holds the headers) and runs down column B.
$firstRow = 2
$lastRow = $ws.Range("B2").End(4).Row
$count = $lastRow - $firstRow + 1

$scores = $ws.Range("B$firstRow`:B$lastRow")

# 50% cutoff: the "High" group is made up of the 578 samples with the
# largest MethylScoreAML value (out of 1147 total), matching the cutoff
# used upstream when the categorical/binary columns were regenerated.
# LARGE(scores, n) gives the value at that cutoff rank, so any score
# >= this threshold is "High".
$highCount = 578
$threshold = $excel.WorksheetFunction.Large($scores, $highCount)

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $score = $ws.Cells.Item($r, 2).Value()
    if ($score -ge $threshold) {
        $ws.Cells.Item($r, 3).Value = "High"
        $ws.Cells.Item($r, 4).Value = 1
    } else {
        $ws.Cells.Item($r, 3).Value = "Low"
        $ws.Cells.Item($r, 4).Value = 0
    }
}
